$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 52.138213
$ws.Range("H2").Value = 156.414639
$ws.Range("I2").Value = 0.2220849502516424
$ws.Range("J2").Value = 0.2220849502516423
$ws.Range("M2").Value = 9.901044000000001
$ws.Range("N2").Value = 29.703132
$ws.Range("O2").Value = 0.3107525783441034
$ws.Range("P2").Value = 0.3107525783441034
$ws.Range("Q2").Value = 516.2227409943721
$ws.Range("R2").Value = 4646.004668949347
$ws.Range("S2").Value = 0.0690134709021198
$ws.Range("T2").Value = 0.0690134709021198

$ws.Range("G3").Value = 52.138213
$ws.Range("H3").Value = 156.414639
$ws.Range("I3").Value = 0.2220849502516424
$ws.Range("J3").Value = 0.2220849502516423
$ws.Range("M3").Value = 7.971374
$ws.Range("N3").Value = 23.914122
$ws.Range("O3").Value = 0.2501882653430435
$ws.Range("P3").Value = 0.2501882653430435
$ws.Range("Q3").Value = 415.613195514662
$ws.Range("R3").Value = 3740.518759631958
$ws.Range("S3").Value = 0.05556304846225452
$ws.Range("T3").Value = 0.05556304846225451

$ws.Range("G4").Value = 52.138213
$ws.Range("H4").Value = 156.414639
$ws.Range("I4").Value = 0.2220849502516424
$ws.Range("J4").Value = 0.2220849502516423
$ws.Range("M4").Value = 3.206223
$ws.Range("N4").Value = 9.618669000000001
$ws.Range("O4").Value = 0.1006300006338893
$ws.Range("P4").Value = 0.1006300006338893
$ws.Range("Q4").Value = 167.166737699499
$ws.Range("R4").Value = 1504.500639295491
$ws.Range("S4").Value = 0.02234840868460005
$ws.Range("T4").Value = 0.02234840868460005

$ws.Range("G5").Value = 52.138213
$ws.Range("H5").Value = 156.414639
$ws.Range("I5").Value = 0.2220849502516424
$ws.Range("J5").Value = 0.2220849502516423
$ws.Range("M5").Value = 10.78286133333333
$ws.Range("N5").Value = 32.348584
$ws.Range("O5").Value = 0.3384291556789638
$ws.Range("P5").Value = 0.3384291556789638
$ws.Range("Q5").Value = 562.1991209467974
$ws.Range("R5").Value = 5059.792088521176
$ws.Range("S5").Value = 0.07516002220266801
$ws.Range("T5").Value = 0.07516002220266799

$ws.Range("G6").Value = 67.324
$ws.Range("H6").Value = 201.972
$ws.Range("I6").Value = 0.2867694600645705
$ws.Range("J6").Value = 0.2867694600645705
$ws.Range("M6").Value = 9.901044000000001
$ws.Range("N6").Value = 29.703132
$ws.Range("O6").Value = 0.3107525783441034
$ws.Range("P6").Value = 0.3107525783441034
$ws.Range("Q6").Value = 666.5778862560001
$ws.Range("R6").Value = 5999.200976304001
$ws.Range("S6").Value = 0.08911434910541167
$ws.Range("T6").Value = 0.08911434910541167

$ws.Range("G7").Value = 67.324
$ws.Range("H7").Value = 201.972
$ws.Range("I7").Value = 0.2867694600645705
$ws.Range("J7").Value = 0.2867694600645705
$ws.Range("M7").Value = 7.971374
$ws.Range("N7").Value = 23.914122
$ws.Range("O7").Value = 0.2501882653430435
$ws.Range("P7").Value = 0.2501882653430435
$ws.Range("Q7").Value = 536.664783176
$ws.Range("R7").Value = 4829.983048584
$ws.Range("S7").Value = 0.07174635376691609
$ws.Range("T7").Value = 0.07174635376691609

$ws.Range("G8").Value = 67.324
$ws.Range("H8").Value = 201.972
$ws.Range("I8").Value = 0.2867694600645705
$ws.Range("J8").Value = 0.2867694600645705
$ws.Range("M8").Value = 3.206223
$ws.Range("N8").Value = 9.618669000000001
$ws.Range("O8").Value = 0.1006300006338893
$ws.Range("P8").Value = 0.1006300006338893
$ws.Range("Q8").Value = 215.855757252
$ws.Range("R8").Value = 1942.701815268
$ws.Range("S8").Value = 0.02885761094807783
$ws.Range("T8").Value = 0.02885761094807783

$ws.Range("G9").Value = 67.324
$ws.Range("H9").Value = 201.972
$ws.Range("I9").Value = 0.2867694600645705
$ws.Range("J9").Value = 0.2867694600645705
$ws.Range("M9").Value = 10.78286133333333
$ws.Range("N9").Value = 32.348584
$ws.Range("O9").Value = 0.3384291556789638
$ws.Range("P9").Value = 0.3384291556789638
$ws.Range("Q9").Value = 725.9453564053334
$ws.Range("R9").Value = 6533.508207648001
$ws.Range("S9").Value = 0.09705114624416493
$ws.Range("T9").Value = 0.0970511462441649

$ws.Range("G10").Value = 64.99978900000001
$ws.Range("H10").Value = 194.999367
$ws.Range("I10").Value = 0.2768693838132169
$ws.Range("J10").Value = 0.2768693838132169
$ws.Range("M10").Value = 9.901044000000001
$ws.Range("N10").Value = 29.703132
$ws.Range("O10").Value = 0.3107525783441034
$ws.Range("P10").Value = 0.3107525783441034
$ws.Range("Q10").Value = 643.5657708797161
$ws.Range("R10").Value = 5792.091937917444
$ws.Range("S10").Value = 0.08603787488450031
$ws.Range("T10").Value = 0.08603787488450031

$ws.Range("G11").Value = 64.99978900000001
$ws.Range("H11").Value = 194.999367
$ws.Range("I11").Value = 0.2768693838132169
$ws.Range("J11").Value = 0.2768693838132169
$ws.Range("M11").Value = 7.971374
$ws.Range("N11").Value = 23.914122
$ws.Range("O11").Value = 0.2501882653430435
$ws.Range("P11").Value = 0.2501882653430435
$ws.Range("Q11").Value = 518.137628040086
$ws.Range("R11").Value = 4663.238652360774
$ws.Range("S11").Value = 0.06926947086282605
$ws.Range("T11").Value = 0.06926947086282605

$ws.Range("G12").Value = 64.99978900000001
$ws.Range("H12").Value = 194.999367
$ws.Range("I12").Value = 0.2768693838132169
$ws.Range("J12").Value = 0.2768693838132169
$ws.Range("M12").Value = 3.206223
$ws.Range("N12").Value = 9.618669000000001
$ws.Range("O12").Value = 0.1006300006338893
$ws.Range("P12").Value = 0.1006300006338893
$ws.Range("Q12").Value = 208.403818486947
$ws.Range("R12").Value = 1875.634366382523
$ws.Range("S12").Value = 0.02786136626862856
$ws.Range("T12").Value = 0.02786136626862856

$ws.Range("G13").Value = 64.99978900000001
$ws.Range("H13").Value = 194.999367
$ws.Range("I13").Value = 0.2768693838132169
$ws.Range("J13").Value = 0.2768693838132169
$ws.Range("M13").Value = 10.78286133333333
$ws.Range("N13").Value = 32.348584
$ws.Range("O13").Value = 0.3384291556789638
$ws.Range("P13").Value = 0.3384291556789638
$ws.Range("Q13").Value = 700.8837114829255
$ws.Range("R13").Value = 6307.953403346329
$ws.Range("S13").Value = 0.09370067179726195
$ws.Range("T13").Value = 0.09370067179726194

$ws.Range("G14").Value = 50.30497766666667
$ws.Range("H14").Value = 150.914933
$ws.Range("I14").Value = 0.2142762058705703
$ws.Range("J14").Value = 0.2142762058705703
$ws.Range("M14").Value = 9.901044000000001
$ws.Range("N14").Value = 29.703132
$ws.Range("O14").Value = 0.3107525783441034
$ws.Range("P14").Value = 0.3107525783441034
$ws.Range("Q14").Value = 498.071797296684
$ws.Range("R14").Value = 4482.646175670156
$ws.Range("S14").Value = 0.06658688345207164
$ws.Range("T14").Value = 0.06658688345207164

$ws.Range("G15").Value = 50.30497766666667
$ws.Range("H15").Value = 150.914933
$ws.Range("I15").Value = 0.2142762058705703
$ws.Range("J15").Value = 0.2142762058705703
$ws.Range("M15").Value = 7.971374
$ws.Range("N15").Value = 23.914122
$ws.Range("O15").Value = 0.2501882653430435
$ws.Range("P15").Value = 0.2501882653430435
$ws.Range("Q15").Value = 400.9997910426473
$ws.Range("R15").Value = 3608.998119383826
$ws.Range("S15").Value = 0.05360939225104687
$ws.Range("T15").Value = 0.05360939225104687

$ws.Range("G16").Value = 50.30497766666667
$ws.Range("H16").Value = 150.914933
$ws.Range("I16").Value = 0.2142762058705703
$ws.Range("J16").Value = 0.2142762058705703
$ws.Range("M16").Value = 3.206223
$ws.Range("N16").Value = 9.618669000000001
$ws.Range("O16").Value = 0.1006300006338893
$ws.Range("P16").Value = 0.1006300006338893
$ws.Range("Q16").Value = 161.288976409353
$ws.Range("R16").Value = 1451.600787684177
$ws.Range("S16").Value = 0.02156261473258289
$ws.Range("T16").Value = 0.02156261473258289

$ws.Range("G17").Value = 50.30497766666667
$ws.Range("H17").Value = 150.914933
$ws.Range("I17").Value = 0.2142762058705703
$ws.Range("J17").Value = 0.2142762058705703
$ws.Range("M17").Value = 10.78286133333333
$ws.Range("N17").Value = 32.348584
$ws.Range("O17").Value = 0.3384291556789638
$ws.Range("P17").Value = 0.3384291556789638
$ws.Range("Q17").Value = 542.4315985560969
$ws.Range("R17").Value = 4881.884387004872
$ws.Range("S17").Value = 0.07251731543486895
$ws.Range("T17").Value = 0.07251731543486893
